# Apply updated DAP (Daily Assessment/Price?) figures to the WESM exposure sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (HOUR 1)
$ws.Range("B2").Value = 28055.51662633527
$ws.Range("C2").Value = 42500
$ws.Range("D2").Value = -14444.48337366473

# Row 3 (HOUR 2)
$ws.Range("B3").Value = 26905.70154390091
$ws.Range("C3").Value = 22500
$ws.Range("D3").Value = 4405.701543900908

# Row 4 (HOUR 3)
$ws.Range("B4").Value = 25723.37458066656
$ws.Range("D4").Value = 3223.374580666561

# Row 5 (HOUR 4)
$ws.Range("B5").Value = 24845.72510195928
$ws.Range("D5").Value = 2345.72510195928

# Row 6 (HOUR 5)
$ws.Range("B6").Value = 24988.90498752365
$ws.Range("C6").Value = 22500
$ws.Range("D6").Value = 2488.904987523649

# Row 7 (HOUR 6)
$ws.Range("B7").Value = 25835.11014617304
$ws.Range("C7").Value = 22500
$ws.Range("D7").Value = 3335.11014617304

# Row 8 (HOUR 7)
$ws.Range("B8").Value = 26401.55040670632
$ws.Range("C8").Value = 22500
$ws.Range("D8").Value = 3901.550406706319

# Row 9 (HOUR 8)
$ws.Range("B9").Value = 28889.35044967622
$ws.Range("C9").Value = 22500
$ws.Range("D9").Value = 6389.350449676222

# Row 10 (HOUR 9)
$ws.Range("B10").Value = 33516.87245290272
$ws.Range("D10").Value = -8983.127547097283

# Row 11 (HOUR 10)
$ws.Range("B11").Value = 29249.825
$ws.Range("D11").Value = -13250.175

# Row 12 (HOUR 11)
$ws.Range("B12").Value = 29168.172
$ws.Range("D12").Value = -23331.828

# Row 13 (HOUR 12)
$ws.Range("B13").Value = 30065.795
$ws.Range("D13").Value = -22434.205

# Row 14 (HOUR 13)
$ws.Range("B14").Value = 30878.246
$ws.Range("D14").Value = -21621.754

# Row 15 (HOUR 14)
$ws.Range("B15").Value = 32676.0655
$ws.Range("D15").Value = -42323.9345

# Row 23 (HOUR 22)
$ws.Range("C23").Value = 74000
$ws.Range("D23").Value = -35714.7935

# Row 24 (HOUR 23)
$ws.Range("C24").Value = 71000

# Row 25 (HOUR 24)
$ws.Range("C25").Value = 67500
